$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block of "driverTitleTestData" test data, mirroring the existing
# emiCalculator block (header / column-labels / values rows).
$ws.Range("A5").Value = "driverTitleTestData"

$ws.Range("A6").Value = "WebPageURL"
$ws.Range("B6").Value = "title"

$ws.Range("A7").Value = "https://www.emicalculator.net"
$ws.Range("B7").Value = "EMI Calculator for Home Loan, Car Loan & Personal Loan in India"

# Wire up the actual hyperlink, then give A7 the same "Hyperlink" look as
# A3 (reuse its cell format so the same style index is shared).
$ws.Hyperlinks.Add($ws.Range("A7"), "https://www.emicalculator.net/") | Out-Null
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# Move the active selection like the authored workbook.
$ws.Range("A9").Select()
